# Added periodic & upfront related scenarios
# The "repaymentstrategy" value (B17) on the ProductLoanInput sheet changes
# from "Mifos style" to "Penalties, Fees, Interest, Principal order", and the
# cell gets a left/top aligned style. Selection is moved to B17 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"
$cell.HorizontalAlignment = -4131
$cell.VerticalAlignment = -4160

$ws.Range("B17").Select()
